$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 412025
$ws.Range("B28").Value = "AWS Core"
$ws.Range("C28").Value = "NoSuchVersion - The version ID specified in the request does not match an existing version."

$ws.Range("A29").Value = 412026
$ws.Range("B29").Value = "AWS Core"
$ws.Range("C29").Value = "NotImplemented - A provided header implies functionality that is not implemented."

$ws.Range("A30").Select()
